$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample")

# Update the LastName / FirstName values for rows 3-5 (Katie Ball, Donald Trump, Mohammed Salah)
$ws.Range("A3").Value = "Katie"
$ws.Range("C3").Value = "Ball"

$ws.Range("A4").Value = "Donald"
$ws.Range("C4").Value = "Trump"

$ws.Range("A5").Value = "Mohammed"
$ws.Range("C5").Value = "Salah"

# Update the selected cell/range on the sheet
$ws.Range("C6").Select()
